# Adds the BiBBS_Geographic sheet with its ch_core_geog / house_nbhd
# data-dictionary rows, and wraps the range in an Excel Table (Table5),
# matching the other two data-dictionary sheets already in the workbook.

$wb = $excel.ActiveWorkbook

# New sheet goes after the last existing sheet (BiBBS_CohortInfo).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "BiBBS_Geographic"

# Column widths (approximate Excel "characters" units for the stored
# OOXML widths used by the sibling sheets).
$ws.Columns.Item(1).ColumnWidth = 16.0
$ws.Columns.Item(2).ColumnWidth = 12.0
$ws.Columns.Item(3).ColumnWidth = 22.0
$ws.Columns.Item(4).ColumnWidth = 52.0
$ws.Columns.Item(5).ColumnWidth = 80.0
$ws.Columns.Item(6).ColumnWidth = 10.0
$ws.Columns.Item(7).ColumnWidth = 11.0
$ws.Columns.Item(8).ColumnWidth = 10.0
$ws.Columns.Item(9).ColumnWidth = 91.0

$headers = @("project", "table", "variable", "full_name", "label", "value_type", "description", "categories", "categories_label")
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

$data = @(
    ,@("BiBBS_Geographic", "ch_core_geog", "date_address_data", "BiBBS_Geographic.ch_core_geog.date_address_data", "Date of latest available address data", "date", $null, $null, $null)
    ,@("BiBBS_Geographic", "ch_core_geog", "age_m", "BiBBS_Geographic.ch_core_geog.age_m", "Participants actual age (months)", "decimal", $null, $null, $null)
    ,@("BiBBS_Geographic", "ch_core_geog", "age_closest_data_point", "BiBBS_Geographic.ch_core_geog.age_closest_data_point", "Age at data point closest to participants actual age", "decimal", $null, $null, $null)
    ,@("BiBBS_Geographic", "ch_core_geog", "temporal_accuracy_m", "BiBBS_Geographic.ch_core_geog.temporal_accuracy_m", "Difference in months between participants age at closest data point and actual a", "decimal", $null, $null, $null)
    ,@("BiBBS_Geographic", "ch_core_geog", "LSOA11CD", "BiBBS_Geographic.ch_core_geog.LSOA11CD", "LSOA 2011 code", "text", $null, $null, $null)
    ,@("BiBBS_Geographic", "ch_core_geog", "WD21CD", "BiBBS_Geographic.ch_core_geog.WD21CD", "Ward 2021 code", "text", $null, $null, $null)
    ,@("BiBBS_Geographic", "ch_core_geog", "is_in_bfd_la", "BiBBS_Geographic.ch_core_geog.is_in_bfd_la", "Is participant in Bradford LA?", "decimal", $null, "0|1", "No|Yes")
    ,@("BiBBS_Geographic", "ch_core_geog", "is_in_bibbs_area", "BiBBS_Geographic.ch_core_geog.is_in_bibbs_area", "Is participant in BiBBS area?", "decimal", $null, "0|1", "No|Yes")
    ,@("BiBBS_Geographic", "ch_core_geog", "data_source", "BiBBS_Geographic.ch_core_geog.data_source", "Source of data: registration (1) or tracing (2)", "decimal", $null, "1|2", "Registration|Tracing")
    ,@("BiBBS_Geographic", "ch_core_geog", "study", "BiBBS_Geographic.ch_core_geog.study", "Is participant in BiB (1) or BiBBS (2)", "decimal", $null, "1|2", "BiB|BiBBS")
    ,@("BiBBS_Geographic", "ch_core_geog", "not_in_eng_wales", "BiBBS_Geographic.ch_core_geog.not_in_eng_wales", "Indicates if address is not in England or Wales", "decimal", $null, "0|1", "No|Yes")
    ,@("BiBBS_Geographic", "ch_core_geog", "missing_address_data", "BiBBS_Geographic.ch_core_geog.missing_address_data", "Indicates if record has missing address data", "decimal", $null, "0|1", "No|Yes")
    ,@("BiBBS_Geographic", "ch_core_geog", "poor_qual_data", "BiBBS_Geographic.ch_core_geog.poor_qual_data", "Indicates if record has poor quality data", "decimal", $null, "0|1", "No|Yes")
    ,@("BiBBS_Geographic", "house_nbhd", "age_m", "BiBBS_Geographic.house_nbhd.age_m", "Participants actual age (months)", "decimal", $null, $null, $null)
    ,@("BiBBS_Geographic", "house_nbhd", "IMD_2010_decile", "BiBBS_Geographic.house_nbhd.IMD_2010_decile", "IMD 2010 decile", "decimal", $null, $null, $null)
    ,@("BiBBS_Geographic", "house_nbhd", "IMD_2010_score", "BiBBS_Geographic.house_nbhd.IMD_2010_score", "IMD 2010 score", "decimal", $null, $null, $null)
    ,@("BiBBS_Geographic", "house_nbhd", "IMD_2019_decile", "BiBBS_Geographic.house_nbhd.IMD_2019_decile", "IMD 2019 decile", "decimal", $null, $null, $null)
    ,@("BiBBS_Geographic", "house_nbhd", "IMD_2019_score", "BiBBS_Geographic.house_nbhd.IMD_2019_score", "IMD 2019 score", "decimal", $null, $null, $null)
    ,@("BiBBS_Geographic", "house_nbhd", "house_type", "BiBBS_Geographic.house_nbhd.house_type", "Property type", "decimal", $null, "1|2|3|4|5", "Terraced|Semi-Detached|Self Contained Flat (Includes Maisonette / Apartment)|Detached|Other")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $value = $row[$c]
        if ($null -ne $value) {
            $ws.Cells.Item($r + 2, $c + 1).Value = $value
        }
    }
}

$lastRow = 20
$lastCol = 9
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table5"
$lo.TableStyle = "TableStyleLight9"

# Match the print setup used by the other data-dictionary sheets.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Restore the original active sheet/tab (sheet creation selects the new sheet).
$wb.Worksheets.Item(1).Activate()

